$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (2 through 7) before writing the new, larger set
$ws.Range("A2:B7").ClearContents()

# Write the brand-new string values first, in the exact order they need to be
# appended to the shared-string table (this mirrors how the source workbook
# was produced), landing directly in their final target cells.
$ws.Range("A4").Value = "BRFS3"
$ws.Range("A5").Value = "CPLE3"
$ws.Range("A6").Value = "ITSA3"
$ws.Range("A7").Value = "USIM3"
$ws.Range("A11").Value = "Tesouro IPCA+ 2035"
$ws.Range("A12").Value = "Tesouro IPCA+ 2045"
$ws.Range("A10").Value = "Tesouro IPCA+ 2026"
$ws.Range("A13").Value = "Tesouro IPCA+ com Juros Semestrais 2035"
$ws.Range("A14").Value = "Tesouro IPCA+ com Juros Semestrais 2055"
$ws.Range("A15").Value = "Tesouro Prefixado 2024"
$ws.Range("A16").Value = "Tesouro Prefixado 2025"
$ws.Range("A17").Value = "Tesouro Prefixado 2026"
$ws.Range("A18").Value = "Tesouro SELIC 2024"
$ws.Range("A19").Value = "Tesouro SELIC 2027"

# Now fill in the remaining cells (reusing already-existing shared strings)
$ws.Range("A2").Value = "BBAS3"
$ws.Range("B2").Value = "Ações"
$ws.Range("A3").Value = "PETR4"
$ws.Range("B3").Value = "Ações"
$ws.Range("B4").Value = "Ações"
$ws.Range("B5").Value = "Ações"
$ws.Range("B6").Value = "Ações"
$ws.Range("B7").Value = "Ações"
$ws.Range("A8").Value = "HASH11"
$ws.Range("B8").Value = "ETF"
$ws.Range("A9").Value = "CPTS11"
$ws.Range("B9").Value = "FII"
$ws.Range("B10").Value = "Tesouro Direto"
$ws.Range("B11").Value = "Tesouro Direto"
$ws.Range("B12").Value = "Tesouro Direto"
$ws.Range("B13").Value = "Tesouro Direto"
$ws.Range("B14").Value = "Tesouro Direto"
$ws.Range("B15").Value = "Tesouro Direto"
$ws.Range("B16").Value = "Tesouro Direto"
$ws.Range("B17").Value = "Tesouro Direto"
$ws.Range("B18").Value = "Tesouro Direto"
$ws.Range("B19").Value = "Tesouro Direto"

# Column A needs to widen to fit the new, longer strings
$ws.Columns.Item(1).ColumnWidth = 37.6
